$d = $word.ActiveDocument

function Find-ParagraphIndexByText($doc, $exactText) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -eq $exactText) {
            return $idx
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. Delete the empty "ListParagraph" bullet sitting between
#    "What type of crime might you be exposed to in certain LLSOAs?" and
#    "Where is the safest place to park your bike?"
# ---------------------------------------------------------------------------
$idxWhatType = Find-ParagraphIndexByText $d "What type of crime might you be exposed to in certain LLSOAs?`r"
$emptyIdx = $idxWhatType + 1
$pEmpty = $d.Paragraphs($emptyIdx)
if ($pEmpty.Range.Text -ne "`r") {
    throw "Unexpected paragraph content when looking for the empty bullet: [$($pEmpty.Range.Text)]"
}
[void]$pEmpty.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Delete the two "D4" / "D5" list paragraphs entirely.
# ---------------------------------------------------------------------------
$idxD4 = Find-ParagraphIndexByText $d "D4`r"
$idxD5 = Find-ParagraphIndexByText $d "D5`r"
$pD4 = $d.Paragraphs($idxD4)
$pD5 = $d.Paragraphs($idxD5)
$delRange = $d.Range($pD4.Range.Start, $pD5.Range.End)
[void]$delRange.Delete()

# ---------------------------------------------------------------------------
# 3. Move <w:lastRenderedPageBreak/> from the "Exploring and cleaning of data"
#    run to the "Allocate tasks to group members" run.
# ---------------------------------------------------------------------------
$idxExploring = Find-ParagraphIndexByText $d "Exploring and cleaning of data`r"
$pExploring = $d.Paragraphs($idxExploring)
$xmlExploring = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="417C82F0" w14:textId="0C35F5C5" w:rsidR="008163A9" w:rsidRDefault="008163A9" w:rsidP="00B93186"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Exploring and cleaning of data</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$pExploring.Range.InsertXML($xmlExploring)

$idxAllocate = Find-ParagraphIndexByText $d "Allocate tasks to group members`r"
$pAllocate = $d.Paragraphs($idxAllocate)
$xmlAllocate = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4447330B" w14:textId="61A919C7" w:rsidR="008163A9" w:rsidRDefault="008163A9" w:rsidP="00B93186"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>Allocate tasks to group members</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$pAllocate.Range.InsertXML($xmlAllocate)

# ---------------------------------------------------------------------------
# 4. Append " University" after "...bike theft Warwick" as its own run.
# ---------------------------------------------------------------------------
$idxWarwick = Find-ParagraphIndexByText $d "Further analysis where we find interesting things – e.g. bike theft Warwick`r"
$pWarwick = $d.Paragraphs($idxWarwick)
$xmlWarwick = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="02D18641" w14:textId="77777777" w:rsidR="00121F3F" w:rsidRDefault="00121F3F" w:rsidP="00121F3F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Further analysis where we find interesting things \u2013 e.g. bike theft Warwick</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> University</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xmlWarwick = $xmlWarwick -replace [regex]::Escape("\u2013"), [char]0x2013
[void]$pWarwick.Range.InsertXML($xmlWarwick)
